$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New cell E3: a quote-prefixed literal text that looks like a formula.
$ws.Range("E3").Value = "'=`$B`$2+`$D`$2 *(COS((A8+`$A`$2)*PI()/180)^2)"

# 2. New rows 31-32: transposed copy of the Angle/Intensity data table.
$ws.Range("A31").Value = "Angle"
$angles = @(-85,-75,-65,-55,-45,-35,-25,-15,-5,5,15,25,35,45,55,65,75,85)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S")
for ($i = 0; $i -lt $angles.Length; $i++) {
    $ws.Range("$($cols[$i])31").Value = $angles[$i]
}

$ws.Range("A32").Value = "Intensity"
$intensities = @(0.475,0.34,0.215,0.125,0.08,0.09,0.15,0.265,0.415,0.525,0.71,0.82,0.87,0.96,0.98,0.88,0.78,0.66)
for ($i = 0; $i -lt $intensities.Length; $i++) {
    $ws.Range("$($cols[$i])32").Value = $intensities[$i]
}

# 3. Row 34: new labels.
$ws.Range("D34").Value = "Enaught"
$ws.Range("E34").Value = "angle"

# 4. Row 35: values used by the calculations below.
$ws.Range("D35").Value = 1.1
$ws.Range("E35").Value = 45

# 5. Row 36: helper formula.
$ws.Range("J36").Formula = "=1.1"

# 6. Row 38 (entered before row 37, matching the shared-string ordering).
$ws.Range("D38").Value = "Itot"
$ws.Range("C38").Value = "pi/2"
$ws.Range("E38").Formula = "=0.5*(D35^2)"
$ws.Range("H38").Formula = "=E38*(COS(90-`$E`$35)^2)"

# 7. Row 37.
$ws.Range("C37").Value = "pi/4"
$ws.Range("D37").Value = "Itot"
$ws.Range("E37").Formula = "=((D35^2)/2)*(1+(SQRT(2)/2)*SIN(2*E35))"
$ws.Range("G37").Value = "Itot"
$ws.Range("H37").Formula = "=E37*(COS(90-`$E`$35)^2)"
$ws.Range("J37").Value = 0.4
$ws.Range("K37").Formula = "=J37/J36"

# 8. Update the sheet view to match the final selection/scroll position
#    (best effort - underlying COM runtime does not persist topLeftCell).
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 2
$ws.Range("E38").Select()
